$d = $word.ActiveDocument

$d.Content.Find.Execute("217÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "252÷4=", 2)
$d.Content.Find.Execute("762÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "710÷4=", 2)
$d.Content.Find.Execute("207÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "662÷3=", 2)
$d.Content.Find.Execute("405÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "355÷2=", 2)
$d.Content.Find.Execute("776÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "303÷4=", 2)
$d.Content.Find.Execute("862÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "581÷4=", 2)
$d.Content.Find.Execute("461÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "395÷5=", 2)
$d.Content.Find.Execute("437÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "770÷2=", 2)
$d.Content.Find.Execute("335÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "335÷4=", 2)
$d.Content.Find.Execute("310÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "842÷9=", 2)
$d.Content.Find.Execute("735÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "592÷7=", 2)
$d.Content.Find.Execute("598÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "925÷7=", 2)
$d.Content.Find.Execute("598÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "157÷5=", 2)
$d.Content.Find.Execute("218÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "841÷2=", 2)
$d.Content.Find.Execute("843÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "464÷9=", 2)
$d.Content.Find.Execute("296÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "185÷2=", 2)
$d.Content.Find.Execute("178÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "180÷9=", 2)
$d.Content.Find.Execute("810÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "977÷5=", 2)
$d.Content.Find.Execute("144÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "239÷8=", 2)
$d.Content.Find.Execute("437÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "653÷6=", 2)
$d.Content.Find.Execute("873÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "294÷4=", 2)
$d.Content.Find.Execute("345÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "136÷2=", 2)
$d.Content.Find.Execute("967÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "251÷9=", 2)
$d.Content.Find.Execute("441÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "896÷8=", 2)
$d.Content.Find.Execute("180÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "152÷9=", 2)
